# Weekly update: a new week's price record is inserted at row 309,
# pushing the existing historical rows (309-382) down by one (310-383).
# The newly freed row 309 is populated with this week's new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(309).Insert()

$ws.Cells.Item(309, 1).Value = 4
$ws.Cells.Item(309, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(309, 3).Value = "Los Lagos"
$ws.Cells.Item(309, 4).Value = 44785
$ws.Cells.Item(309, 5).Value = 10
$ws.Cells.Item(309, 6).Value = 100112023
$ws.Cells.Item(309, 7).Value = "Brócoli"
$ws.Cells.Item(309, 8).Value = "Sin especificar"
$ws.Cells.Item(309, 9).Value = "Primera"
$ws.Cells.Item(309, 10).Value = 1500
$ws.Cells.Item(309, 11).Value = 1500
$ws.Cells.Item(309, 12).Value = 1500
$ws.Cells.Item(309, 13).Value = 1500
$ws.Cells.Item(309, 14).Value = "$/unidad"
$ws.Cells.Item(309, 15).Value = "Región Metropolitana"
$ws.Cells.Item(309, 16).Value = 1500
$ws.Cells.Item(309, 17).Value = 1
$ws.Cells.Item(309, 18).Value = "Hortaliza"
